$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.358.55"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.443.70"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "411.09"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.92"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -3.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +6.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.757"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +12.39%  "
$ws.Range("E10").Value = "  +18.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.38"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.93"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +7.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.54"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +4.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000195"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +53.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.444.09"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Uniswap"
$ws.Range("B17").NumberFormat = "General"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C17").NumberFormat = "General"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.84"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +7.79%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Polygon"
$ws.Range("B18").NumberFormat = "General"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C18").NumberFormat = "General"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.05"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +3.77%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("B19").NumberFormat = "General"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C19").NumberFormat = "General"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "62.298.70"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "406.77"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +30.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "89.50"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +7.04%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.37"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +5.14%  "
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "32.19"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +9.27%  "
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.53"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +3.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.75"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +10.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "44.16"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +7.29%  "
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("E33").Value = "  +4.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +2.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "52.40"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.40"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  +7.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.314"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +8.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.59"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.04"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.81"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("E46").Value = "  +4.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.14"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.124.49"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0372"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +8.33%  "
